# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    3 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    4 = @{ B = 0.6545652718822623; C = 0.04103571897497393; D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.950081467445961 }
    5 = @{ B = 0.1169995834814548; C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987;  G = 2.426980108624251 }
    6 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 13.86384647080068;   G = 19.48425592650926 }
    7 = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 3.223369029078222;  E = 13.86384647080068;   G = 20.15985084044064 }
    8 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
